# Updated cryptos list - apply latest price/volume(1h) scrape values.
# Note: several "Price" column values look like plain decimal numbers
# (e.g. 590.93) but must remain text, matching the source data which mixes
# thousand-separated prices ("2.528.44") with plain decimals. A leading
# apostrophe forces Excel to keep them as text instead of auto-converting
# to a floating point number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.257.31"
$ws.Range("E2").Value = "  +1.26%  "
$ws.Range("D3").Value = "2.528.44"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'590.93"
$ws.Range("E5").Value = "  +1.38%  "
$ws.Range("D6").Value = "'173.89"
$ws.Range("E6").Value = "  +4.67%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +0.21%  "
$ws.Range("D9").Value = "2.527.64"
$ws.Range("E9").Value = "  -2.03%  "
$ws.Range("D10").Value = "'0.139"
$ws.Range("E10").Value = "  +1.28%  "
$ws.Range("E11").Value = "  +2.12%  "
$ws.Range("E12").Value = "  +0.01%  "
$ws.Range("E13").Value = "  -3.24%  "
$ws.Range("E14").Value = "  -0.10%  "
$ws.Range("D15").Value = "2.989.12"
$ws.Range("E15").Value = "  -1.94%  "
$ws.Range("E16").Value = "  -0.13%  "
$ws.Range("D17").Value = "67.024.88"
$ws.Range("E17").Value = "  +1.10%  "
$ws.Range("D18").Value = "2.530.11"
$ws.Range("E18").Value = "  -1.94%  "
$ws.Range("D19").Value = "'8.07"
$ws.Range("E19").Value = "  +4.97%  "
$ws.Range("D20").Value = "'11.39"
$ws.Range("E20").Value = "  +0.07%  "
$ws.Range("D21").Value = "'354.94"
$ws.Range("E21").Value = "  +1.17%  "
$ws.Range("E22").Value = "  -0.58%  "
$ws.Range("D23").Value = "'4.61"
$ws.Range("E23").Value = "  +0.69%  "
$ws.Range("E24").Value = "  +6.86%  "
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("E26").Value = "  +1.69%  "
$ws.Range("D27").Value = "'9.92"
$ws.Range("E27").Value = "  +0.20%  "
$ws.Range("E28").Value = "  -0.34%  "
$ws.Range("D29").Value = "2.655.86"
$ws.Range("E29").Value = "  -2.13%  "
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("D31").Value = "'539.72"
$ws.Range("E31").Value = "  +2.02%  "
$ws.Range("D32").Value = "'8.18"
$ws.Range("E32").Value = "  +2.22%  "
$ws.Range("E33").Value = "  +1.08%  "
$ws.Range("E34").Value = "  +0.37%  "
$ws.Range("D35").Value = "'0.130"
$ws.Range("E35").Value = "  -0.93%  "
$ws.Range("D36").Value = "'0.999"
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("E37").Value = "  +0.74%  "
$ws.Range("D38").Value = "'157.48"
$ws.Range("E38").Value = "  +0.60%  "
$ws.Range("E39").Value = "  -0.26%  "
$ws.Range("D40").Value = "'18.44"
$ws.Range("E40").Value = "  +0.92%  "
$ws.Range("D41").Value = "'0.354"
$ws.Range("E41").Value = "  -1.55%  "
$ws.Range("E42").Value = "  +1.55%  "
$ws.Range("E43").Value = "  +1.20%  "
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("D45").Value = "'2.51"
$ws.Range("E45").Value = "  +4.22%  "
$ws.Range("D46").Value = "'149.09"
$ws.Range("E46").Value = "  +0.33%  "
$ws.Range("E47").Value = "  -1.09%  "
$ws.Range("D48").Value = "0.0₆0277"
$ws.Range("E48").Value = "  -2.68%  "
$ws.Range("E49").Value = "  -0.11%  "
$ws.Range("D50").Value = "'1.69"
$ws.Range("E50").Value = "  -0.35%  "
$ws.Range("E51").Value = "  -0.18%  "
